$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list data (prices in column D, volume % in column E,
# plus two coin-name/link row swaps at rows 28/29 and 39/40).
# Column D values are prefixed with a leading apostrophe to force Excel
# to store them as text (matching the original inlineStr cell type and
# preserving exact formatting such as trailing zeros / multi-dot prices).

$ws.Range("D2").Value = "'86.608.45"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "'3.167.64"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'202.64"
$ws.Range("E5").Value = "  -7.68%  "
$ws.Range("D6").Value = "'602.28"
$ws.Range("E6").Value = "  -7.53%  "
$ws.Range("E7").Value = "  -8.55%  "
$ws.Range("D8").Value = "'0.653"
$ws.Range("E8").Value = "  +7.62%  "
$ws.Range("D9").Value = "'0.998"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "'3.169.55"
$ws.Range("E10").Value = "  -5.00%  "
$ws.Range("D11").Value = "'0.527"
$ws.Range("E11").Value = "  -10.67%  "
$ws.Range("D12").Value = "'0.176"
$ws.Range("E12").Value = "  +4.39%  "
$ws.Range("E13").Value = "  -15.67%  "
$ws.Range("D14").Value = "'3.748.30"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("D15").Value = "'5.19"
$ws.Range("E15").Value = "  -5.12%  "
$ws.Range("D16").Value = "'86.611.25"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "'31.66"
$ws.Range("E17").Value = "  -10.14%  "
$ws.Range("D18").Value = "'3.184.66"
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").Value = "'13.27"
$ws.Range("E19").Value = "  -9.42%  "
$ws.Range("D20").Value = "'2.87"
$ws.Range("E20").Value = "  -8.97%  "
$ws.Range("D21").Value = "'409.58"
$ws.Range("E21").Value = "  -10.43%  "
$ws.Range("D22").Value = "'8.35"
$ws.Range("E22").Value = "  -14.49%  "
$ws.Range("D23").Value = "'4.96"
$ws.Range("E23").Value = "  -9.88%  "
$ws.Range("D24").Value = "'5.01"
$ws.Range("E24").Value = "  -10.56%  "
$ws.Range("D25").Value = "'11.32"
$ws.Range("E25").Value = "  -11.61%  "
$ws.Range("D26").Value = "'3.350.94"
$ws.Range("E26").Value = "  -4.68%  "
$ws.Range("D27").Value = "'72.52"
$ws.Range("E27").Value = "  -7.71%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0000127"
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("D30").Value = "'0.162"
$ws.Range("E30").Value = "  -18.77%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "'528.66"
$ws.Range("E32").Value = "  -10.99%  "
$ws.Range("D33").Value = "'8.16"
$ws.Range("E33").Value = "  -13.21%  "
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  -14.33%  "
$ws.Range("E35").Value = "  -22.71%  "
$ws.Range("D36").Value = "'6.43"
$ws.Range("E36").Value = "  -11.59%  "
$ws.Range("D37").Value = "'0.131"
$ws.Range("E37").Value = "  -10.40%  "
$ws.Range("D38").Value = "'21.48"
$ws.Range("E38").Value = "  -8.12%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'21.80"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'2.90"
$ws.Range("E42").Value = "  -6.72%  "
$ws.Range("D43").Value = "'0.369"
$ws.Range("E43").Value = "  -12.53%  "
$ws.Range("E44").Value = "  -15.35%  "
$ws.Range("D45").Value = "'146.70"
$ws.Range("E45").Value = "  -7.73%  "
$ws.Range("D46").Value = "'170.59"
$ws.Range("E46").Value = "  -10.48%  "
$ws.Range("D47").Value = "'42.90"
$ws.Range("E47").Value = "  -7.85%  "
$ws.Range("D48").Value = "'0.125"
$ws.Range("E48").Value = "  +8.98%  "
$ws.Range("E49").Value = "  -14.07%  "
$ws.Range("D50").Value = "'3.90"
$ws.Range("E50").Value = "  -11.84%  "
$ws.Range("D51").Value = "'0.582"
$ws.Range("E51").Value = "  -12.33%  "
